$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row-level updates: Coin, Link (text), Price, Volume(1h), Hora (forced text)
$rowUpdates = @{
    2 = @{ D="306.81"; E="-3.33%"; G="17" }
    3 = @{ D="40.13"; E="-4.35%"; G="17" }
    4 = @{ D="5.042"; E="-2.72%"; G="17" }
    5 = @{ D="0.07600"; E="-6.33%"; G="17" }
    6 = @{ D="4.255"; E="-2.67%"; G="17" }
    7 = @{ D="1.589"; E="-8.99%"; G="17" }
    8 = @{ D="0.9087"; E="-2.23%"; G="17" }
    9 = @{ E="-10.45%"; G="17" }
    10 = @{ D="0.1755"; E="-5.70%"; G="17" }
    11 = @{ D="0.08988"; E="-2.32%"; G="17" }
    12 = @{ D="0.04366"; E="-5.04%"; G="17" }
    13 = @{ D="0.1053"; E="-0.01%"; G="17" }
    14 = @{ D="0.001236"; E="-3.08%"; G="17" }
    15 = @{ D="0.005824"; E="2.34%"; G="17" }
    16 = @{ E="2,404.66%"; G="17" }
    17 = @{ D="3.368"; E="0.52%"; G="17" }
    18 = @{ E="-3.34%"; G="17" }
    19 = @{ D="0.3280"; E="-3.18%"; G="17" }
    20 = @{ D="6.887"; E="-6.50%"; G="17" }
    21 = @{ D="0.1351"; E="-3.06%"; G="17" }
    22 = @{ D="0.2821"; E="8.29%"; G="17" }
    23 = @{ D="0.04160"; E="-0.48%"; G="17" }
    24 = @{ D="0.001225"; E="-1.74%"; G="17" }
    25 = @{ D="0.004062"; E="-4.73%"; G="17" }
    26 = @{ D="0.0001304"; E="6.48%"; G="17" }
    27 = @{ G="17" }
    28 = @{ G="17" }
    29 = @{ G="17" }
    30 = @{ G="17" }
    31 = @{ G="17" }
    32 = @{ G="17" }
    33 = @{ G="17" }
    34 = @{ G="17" }
    35 = @{ G="17" }
    36 = @{ G="17" }
    37 = @{ G="17" }
    38 = @{ D="0.02414"; E="-6.43%"; G="17" }
    39 = @{ D="0.05129"; E="-6.15%"; G="17" }
    40 = @{ D="0.007858"; E="-3.32%"; G="17" }
    41 = @{ D="0.1307"; E="-6.01%"; G="17" }
    42 = @{ D="0.007094"; E="8.60%"; G="17" }
    43 = @{ D="0.001974"; E="-5.79%"; G="17" }
    44 = @{ D="0.008362"; E="9.32%"; G="17" }
    45 = @{ D="0.3320"; E="-4.12%"; G="17" }
    46 = @{ D="0.00006466"; E="-4.44%"; G="17" }
    47 = @{ D="0.00000000752"; E="-0.08%"; G="17" }
    48 = @{ B="CoinbaseStockToken"; C="https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"; D="0.003003"; E="-27.04%"; G="17" }
    49 = @{ B="BOLO"; C="https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"; D="0.004877"; E="44.17%"; G="17" }
    50 = @{ D="0.00002107"; E="-0.08%"; G="17" }
    51 = @{ D="0.0002007"; E="-0.08%"; G="17" }
}

# Columns whose values must stay as TEXT even though they look numeric
$textForcedCols = @("D", "E", "G")

foreach ($row in $rowUpdates.Keys) {
    $vals = $rowUpdates[$row]
    foreach ($col in $vals.Keys) {
        $addr = "$col$row"
        if ($textForcedCols -contains $col) {
            $ws.Range($addr).NumberFormat = "@"
        }
        $ws.Range($addr).Value = $vals[$col]
    }
}
